$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.9212345884845585
$ws.Range("D2").Value = 0.9396192410111824
$ws.Range("C3").Value = 0.8549134131291252
$ws.Range("D3").Value = 0.8871032311863464
$ws.Range("C4").Value = 0.7990294877625025
$ws.Range("D4").Value = 0.8414360929920358
$ws.Range("C5").Value = 0.751101105728752
$ws.Range("D5").Value = 0.8008950318430346
$ws.Range("C6").Value = 0.7096389690341026
$ws.Range("D6").Value = 0.765615050747926
$ws.Range("C7").Value = 0.6742722234351348
$ws.Range("D7").Value = 0.7345111724812472
$ws.Range("C8").Value = 0.6441603250092323
$ws.Range("D8").Value = 0.7069676128769365
$ws.Range("C9").Value = 0.6173755577065544
$ws.Range("D9").Value = 0.6825478613802501
$ws.Range("C10").Value = 0.594064940237079
$ws.Range("D10").Value = 0.6608175274174485
$ws.Range("C11").Value = 0.5738936188979189
$ws.Range("D11").Value = 0.6412901212078319
$ws.Range("C12").Value = 0.5563045006114613
$ws.Range("D12").Value = 0.6238692228113666
$ws.Range("C13").Value = 0.5402558571151572
$ws.Range("D13").Value = 0.6082492485478933
$ws.Range("C14").Value = 0.5259588038490605
$ws.Range("D14").Value = 0.594540446454448
$ws.Range("C15").Value = 0.5133217673152679
$ws.Range("D15").Value = 0.5823060403406269
$ws.Range("C16").Value = 0.5020451817775909
$ws.Range("D16").Value = 0.571686615128736
$ws.Range("C17").Value = 0.4919841815731267
$ws.Range("D17").Value = 0.5612001064250544
$ws.Range("C18").Value = 0.4829983828539188
$ws.Range("D18").Value = 0.552513452638021
$ws.Range("C19").Value = 0.4751468231107248
$ws.Range("D19").Value = 0.5440830707804111
$ws.Range("C20").Value = 0.4679479883834771
$ws.Range("D20").Value = 0.5361897048972942
$ws.Range("C21").Value = 0.4612255226185729
$ws.Range("D21").Value = 0.5290331575418299
$ws.Range("C22").Value = 0.4553711579725309
$ws.Range("D22").Value = 0.5231525247108409
$ws.Range("C23").Value = 0.4498961147230905
$ws.Range("D23").Value = 0.517150432035412
$ws.Range("C24").Value = 0.4446892358286524
$ws.Range("D24").Value = 0.5121802040685866
$ws.Range("C25").Value = 0.4403620223487189
$ws.Range("D25").Value = 0.5077674667663523
$ws.Range("C26").Value = 0.4364558409990663
$ws.Range("D26").Value = 0.5039896771737485
$ws.Range("C27").Value = 0.432978782609613
$ws.Range("D27").Value = 0.5006899482954685
$ws.Range("C28").Value = 0.4294692694028511
$ws.Range("D28").Value = 0.4977483347212673
$ws.Range("C29").Value = 0.4261965955413958
$ws.Range("D29").Value = 0.4950858502496803
$ws.Range("C30").Value = 0.4231808931649826
$ws.Range("D30").Value = 0.4924522993348402
$ws.Range("C31").Value = 0.4207121372285994
$ws.Range("D31").Value = 0.4900050534803722
$ws.Range("C32").Value = 0.4183792934127255
$ws.Range("D32").Value = 0.4879843392316656
$ws.Range("C33").Value = 0.4162970144845211
$ws.Range("D33").Value = 0.4858690914358637
$ws.Range("C34").Value = 0.414401627889147
$ws.Range("D34").Value = 0.4837435210256819
$ws.Range("C35").Value = 0.4126274588216012
$ws.Range("D35").Value = 0.4817945978053172
$ws.Range("C36").Value = 0.4111033678057102
$ws.Range("D36").Value = 0.4801152138556302
$ws.Range("C37").Value = 0.4094160364198506
$ws.Range("D37").Value = 0.4787378436759993
$ws.Range("C38").Value = 0.4080507175216367
$ws.Range("D38").Value = 0.477481684762909
$ws.Range("C39").Value = 0.4069754296545577
$ws.Range("D39").Value = 0.4763357995341209
$ws.Range("C40").Value = 0.4057103435786321
$ws.Range("D40").Value = 0.4752902814153452
$ws.Range("C41").Value = 0.4046661606852489
$ws.Range("D41").Value = 0.4743361505733885
$ws.Range("C42").Value = 0.403787985689256
$ws.Range("D42").Value = 0.4733670593668449
$ws.Range("C43").Value = 0.4031226879127355
$ws.Range("D43").Value = 0.4726809906673672
$ws.Range("C44").Value = 0.3997246006133469
$ws.Range("D44").Value = 0.4693257805294755
$ws.Range("C45").Value = 0.3966755071281107
$ws.Range("D45").Value = 0.4670267892819492
$ws.Range("C46").Value = 0.3938173180274999
$ws.Range("D46").Value = 0.4644368814612637
